# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorganizes the "Estado de Cuenta" detail table (rows 16-29) so the two
# workers' late-payment periods (2106-2112) are interleaved by period
# (Gustavo, Luisa, Gustavo, Luisa, ...) instead of grouped by worker, and
# refreshes the "Salario Basico" (base salary) values used for the newly
# added part of the account statement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tipo Doc(B), N Doc(C), Nombre(D), Periodo(E), Valor Mora(F), Salario Basico(G)
$rows = @(
    @(16, "CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2106", 40000,  1000000),
    @(17, "CC", "37863609", "LUISA EDILMA REY PABON",        "2106", 60000,  1500000),
    @(18, "CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2107", 57380,  1000000),
    @(19, "CC", "37863609", "LUISA EDILMA REY PABON",        "2107", 60000,  1500000),
    @(20, "CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2108", 57380,  1000000),
    @(21, "CC", "37863609", "LUISA EDILMA REY PABON",        "2108", 60000,  1500000),
    @(22, "CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2109", 57380,  1000000),
    @(23, "CC", "37863609", "LUISA EDILMA REY PABON",        "2109", 60000,  1500000),
    @(24, "CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2110", 57380,  1000000),
    @(25, "CC", "37863609", "LUISA EDILMA REY PABON",        "2110", 60000,  1500000),
    @(26, "CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2111", 40000,  1000000),
    @(27, "CC", "37863609", "LUISA EDILMA REY PABON",        "2111", 60000,  1500000),
    @(28, "CC", "73201303", "GUSTAVO ADOLFO GALEANO OVIEDO", "2112", 21333,  1000000),
    @(29, "CC", "37863609", "LUISA EDILMA REY PABON",        "2112", 32000,  1500000)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = $r[1]
    $ws.Range("C$rowNum").Value = $r[2]
    $ws.Range("D$rowNum").Value = $r[3]
    $ws.Range("E$rowNum").Value = $r[4]
    $ws.Range("F$rowNum").Value = $r[5]
    $ws.Range("G$rowNum").Value = $r[6]
}
